$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before the old "Amount reimbursed" column (E).
# This shifts the old E (Amount reimbursed) -> F and old F (Amount left to
# reimburse, with its formula) -> G, and Excel auto-adjusts the formula
# references (E2->F2 etc.) the same way the original diff shows.
$ws.Columns.Item(5).Insert()

# New column E width (best achievable value through this host's rounding;
# the target stored width is 12.42578125, the nearest value reachable
# through the ColumnWidth setter's internal pixel rounding is 12.5).
$ws.Columns.Item(5).ColumnWidth = 11.666666666666666

# New header for column E.
$ws.Range("E1").Value = "Method"

# Payment method for the existing 4 participants.
$ws.Range("E2").Value = "cash"
$ws.Range("E3").Value = "cash"
$ws.Range("E4").Value = "cash"
$ws.Range("E5").Value = "momo"

# New participant: Songsare Amdji Pierre.
$ws.Range("A6").Value = 5
$ws.Range("B6").Value = "Songsare Amdji Pierre"
$ws.Range("C6").Value = 694865719
$ws.Range("D6").Value = 300
$ws.Range("E6").Value = "cash"
$ws.Range("F6").Value = 0

# G4:G6 becomes a shared formula group (matches Excel's behaviour when a
# formula is entered once and filled down over several cells at once).
$ws.Range("G4:G6").Formula = "=D4-F4-300"

# Restore the selection to match the post-edit state recorded in the diff.
$ws.Range("G7").Select() | Out-Null
